# feat: readme and categories
# Build the "Z-test by category" worksheet: a bold/bordered/centered header
# row followed by 12 rows of per-category Z-test results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Категория", "Кол-во сотрудников", "Средняя разница", "Z-статистика", "P-значение", "Вывод")

$data = @(
    @("Мерч и сувениры КРОК", 589, -221.01, -7.337, 0,     "❌ Отрицательный рост"),
    @("Без категории",        671, 296.79,  15.45,  0,     "✅ Рост"),
    @("В стиле КРОК",         372, -59.81,  -1.845, 0.065, "❌ Нет роста"),
    @("Сервис",               232, -130.75, -2.761, 0.006, "❌ Отрицательный рост"),
    @("Лучше вместе",         256, -54.16,  -1.091, 0.275, "❌ Нет роста"),
    @("Поймай баланс",        359, -238.54, -7.766, 0,     "❌ Отрицательный рост"),
    @("КРОК FEST TOUR",       45,  -603.78, -3.281, 0.001, "❌ Отрицательный рост"),
    @("Развивай себя",        50,  -39.9,   -0.39,  0.697, "❌ Нет роста"),
    @("Прокачай эффективность", 180, -105.28, -2.008, 0.045, "❌ Отрицательный рост"),
    @("Смена деятельности",   99,  -12.37,  -0.2,   0.841, "❌ Нет роста"),
    @("Своими руками",        45,  -11.89,  -0.075, 0.94,  "❌ Нет роста"),
    @("Новогодние атрибуты",  41,  281.71,  3.25,   0.001, "✅ Рост")
)

# Header row
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Header styling: bold font, thin box border, centered horizontally, top vertically.
# Apply to A1 first, then propagate the exact same resulting format to the rest
# of the header row via a format-only paste, so every header cell shares one style.
$a1 = $ws.Range("A1")
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1

$a1.Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)

Write-Output "done"
